# Updated capital structure database
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value, applied identically to rows 2 and 3
$updates = @{
    "D" = 0.118
    "E" = 1.714
    "K" = 1.95
    "L" = 0.03271812080536913
    "U" = 223
    "V" = 1.275014293882218
    "W" = 0.01058631921824104
    "X" = 0.07734787539784076
    "Y" = -0.06676155617959971
    "Z" = 0.4084989718985607
    "AA" = 0
    "AB" = 0.06276971729480124
    "AC" = -0.06276971729480124
    "AD" = 116.1
    "AF" = 116.1
    "AG" = -106.9
    "AH" = 0.3989690721649484
    "AI" = 0.3859707446808511
    "AJ" = -1.572058823529412
    "AK" = -1.374035989717224
}

foreach ($row in 2,3) {
    foreach ($col in $updates.Keys) {
        $ws.Range("$col$row").Value = $updates[$col]
    }
}
